$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_related_resources")

# Append the new related-resource row (row 3) with the modified data.
$ws.Range("A3").Value = "about-nci/organization/sharpless-nci-director"
$ws.Range("B3").Value = "Press Release"
$ws.Range("C3").Value = "English"

# Make this sheet the active one, with the same selection Excel left behind.
$ws.Activate()
$ws.Range("E19").Select()
